$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.477.63"
$ws.Range("E2").Value = "'  +5.29%  "
$ws.Range("D3").Value = "'1.724.89"
$ws.Range("E3").Value = "'  +4.66%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'225.29"
$ws.Range("E5").Value = "'  +3.10%  "
$ws.Range("D6").Value = "'0.5348"
$ws.Range("E6").Value = "'  +2.90%  "
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("D8").Value = "'0.2659"
$ws.Range("E8").Value = "'  +1.55%  "
$ws.Range("D9").Value = "'0.06586"
$ws.Range("E9").Value = "'  +4.53%  "
$ws.Range("D10").Value = "'21.52"
$ws.Range("E10").Value = "'  +5.94%  "
$ws.Range("D11").Value = "'0.07667"
$ws.Range("E11").Value = "'  +0.07%  "
$ws.Range("D12").Value = "'4.596"
$ws.Range("E12").Value = "'  +0.44%  "
$ws.Range("D13").Value = "'1.727.24"
$ws.Range("E13").Value = "'  +5.19%  "
$ws.Range("D14").Value = "'1.962.49"
$ws.Range("E14").Value = "'  +4.72%  "
$ws.Range("D15").Value = "'0.5784"
$ws.Range("E15").Value = "'  +3.81%  "
$ws.Range("E16").Value = "'  +2.12%  "
$ws.Range("D17").Value = "'67.75"
$ws.Range("E17").Value = "'  +4.17%  "
$ws.Range("D18").Value = "'27.483.96"
$ws.Range("E18").Value = "'  +5.54%  "
$ws.Range("D19").Value = "'218.29"
$ws.Range("E19").Value = "'  +12.44%  "
$ws.Range("E20").Value = "'  +0.04%  "
$ws.Range("D21").Value = "'4.722"
$ws.Range("E21").Value = "'  +2.80%  "
$ws.Range("E22").Value = "'  +1.14%  "
$ws.Range("D23").Value = "'6.027"
$ws.Range("E23").Value = "'  +1.95%  "
$ws.Range("E24").Value = "'  +0.12%  "
$ws.Range("D25").Value = "'142.77"
$ws.Range("E25").Value = "'  -1.08%  "
$ws.Range("D26").Value = "'1.752"
$ws.Range("E26").Value = "'  +15.63%  "
$ws.Range("E27").Value = "'  +4.41%  "
$ws.Range("D28").Value = "'7.322"
$ws.Range("E28").Value = "'  +1.99%  "
$ws.Range("D29").Value = "'16.49"
$ws.Range("E29").Value = "'  +4.08%  "
$ws.Range("D30").Value = "'0.05481"
$ws.Range("E30").Value = "'  +1.48%  "
$ws.Range("D31").Value = "'1.299"
$ws.Range("E31").Value = "'  +2.23%  "
$ws.Range("D32").Value = "'3.554"
$ws.Range("D33").Value = "'3.435"
$ws.Range("E33").Value = "'  +3.40%  "
$ws.Range("D34").Value = "'1.656"
$ws.Range("E34").Value = "'  +6.47%  "
$ws.Range("D35").Value = "'2.858"
$ws.Range("E35").Value = "'  +2.85%  "
$ws.Range("D36").Value = "'0.9559"
$ws.Range("E36").Value = "'  +1.65%  "
$ws.Range("D37").Value = "'2.421"
$ws.Range("E37").Value = "'  +0.29%  "
$ws.Range("D38").Value = "'0.5933"
$ws.Range("E38").Value = "'  +6.46%  "
$ws.Range("D39").Value = "'0.01647"
$ws.Range("E39").Value = "'  +4.72%  "
$ws.Range("D40").Value = "'5.896"
$ws.Range("E40").Value = "'  +2.76%  "
$ws.Range("B41").Value = "'Maker"
$ws.Range("C41").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.048.48"
$ws.Range("E41").Value = "'  +2.13%  "
$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8468"
$ws.Range("E42").Value = "'  +2.95%  "
$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("D44").Value = "'101.32"
$ws.Range("E44").Value = "'  +0.62%  "
$ws.Range("D45").Value = "'1.867.63"
$ws.Range("E45").Value = "'  +4.69%  "
$ws.Range("D46").Value = "'0.0₈118"
$ws.Range("E46").Value = "'  +4.43%  "
$ws.Range("D47").Value = "'58.69"
$ws.Range("E47").Value = "'  +2.57%  "
$ws.Range("D48").Value = "'0.4479"
$ws.Range("E48").Value = "'  +3.76%  "
$ws.Range("D49").Value = "'8.173"
$ws.Range("E49").Value = "'  +3.44%  "
$ws.Range("D50").Value = "'1.002"
$ws.Range("E50").Value = "'  +0.38%  "
$ws.Range("D51").Value = "'0.05250"
$ws.Range("E51").Value = "'  +3.13%  "
